$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1274.2142
$ws.Range("I32").Value = 822.25
$ws.Range("J32").Value = 1455
$ws.Range("K32").Value = 822.25
$ws.Range("L32").Value = 1455
$ws.Range("M32").Value = -496.25
$ws.Range("N32").Value = -2107

$ws.Range("H38").Value = 1112752.2
$ws.Range("I38").Value = 63.4
$ws.Range("J38").Value = 2503613.2
$ws.Range("K38").Value = 190.2
$ws.Range("L38").Value = 7510839.600000001
$ws.Range("M38").Value = 181.8
$ws.Range("N38").Value = -7511583.600000001

$ws.Range("H40").Value = 1704.174
$ws.Range("I40").Value = 1336.25
$ws.Range("J40").Value = 1900.4
$ws.Range("K40").Value = 1336.25
$ws.Range("L40").Value = 1900.4
$ws.Range("M40").Value = -1161.25
$ws.Range("N40").Value = -2250.4

$ws.Range("H113").Value = 3184.7083
$ws.Range("I113").Value = 2996.0588
$ws.Range("J113").Value = 3642.8572
$ws.Range("K113").Value = 2996.0588
$ws.Range("L113").Value = 3642.8572
$ws.Range("M113").Value = 257.9412000000002
$ws.Range("N113").Value = -10150.8572

$ws.Range("H116").Value = 2305.2666
$ws.Range("I116").Value = 2236.3635
$ws.Range("J116").Value = 2494.75
$ws.Range("K116").Value = 2236.3635
$ws.Range("L116").Value = 2494.75
$ws.Range("M116").Value = 1205.6365
$ws.Range("N116").Value = -9378.75

$ws.Range("H129").Value = 2153.6365
$ws.Range("I129").Value = 700
$ws.Range("J129").Value = 2299
$ws.Range("K129").Value = 2100
$ws.Range("L129").Value = 6897
$ws.Range("M129").Value = 2900
$ws.Range("N129").Value = -16897

$ws.Range("H132").Value = 8934087
$ws.Range("I132").Value = 5323.3477
$ws.Range("J132").Value = 50006400
$ws.Range("K132").Value = 15970.0431
$ws.Range("L132").Value = 150019200
$ws.Range("M132").Value = -13440.0431
$ws.Range("N132").Value = -150024260

$ws.Range("H138").Value = 8338500.5
$ws.Range("I138").Value = 20837276
$ws.Range("J138").Value = 5983.9165
$ws.Range("K138").Value = 62511828
$ws.Range("L138").Value = 17951.7495
$ws.Range("M138").Value = -62506688
$ws.Range("N138").Value = -28231.7495

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 25000
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws.Range("H74").Value = 2549.0188
$ws.Range("I74").Value = 629.1372699999999
$ws.Range("J74").Value = 51506
$ws.Range("K74").Value = 629.1372699999999
$ws.Range("L74").Value = 51506
$ws.Range("M74").Value = 244.8627300000001
$ws.Range("N74").Value = -53254

$ws.Range("H77").Value = 2549.0188
$ws.Range("I77").Value = 629.1372699999999
$ws.Range("J77").Value = 51506
$ws.Range("K77").Value = 3145.68635
$ws.Range("L77").Value = 257530
$ws.Range("M77").Value = 1222.31365
$ws.Range("N77").Value = -266266

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1525.25
$ws.Range("I107").Value = 867
$ws.Range("J107").Value = 3500
$ws.Range("K107").Value = 867
$ws.Range("L107").Value = 3500
$ws.Range("M107").Value = 1053
$ws.Range("N107").Value = -7340

$ws.Range("H122").Value = 59408.89
$ws.Range("J122").Value = 59408.89
$ws.Range("L122").Value = 59408.89
$ws.Range("N122").Value = -69208.89

$ws.Range("H133").Value = 46466.668
$ws.Range("J133").Value = 46466.668
$ws.Range("L133").Value = 46466.668
$ws.Range("N133").Value = -56586.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6481.2354
$ws.Range("I31").Value = 2139.5
$ws.Range("J31").Value = 16901.4
$ws.Range("K31").Value = 2139.5
$ws.Range("L31").Value = 16901.4
$ws.Range("M31").Value = -1844.5
$ws.Range("N31").Value = -17491.4

$ws.Range("H34").Value = 6481.2354
$ws.Range("I34").Value = 2139.5
$ws.Range("J34").Value = 16901.4
$ws.Range("K34").Value = 2139.5
$ws.Range("L34").Value = 16901.4
$ws.Range("M34").Value = -1937.5
$ws.Range("N34").Value = -17305.4

$ws.Range("H99").Value = 1692.5834
$ws.Range("I99").Value = 1495.6666
$ws.Range("J99").Value = 2283.3333
$ws.Range("K99").Value = 1495.6666
$ws.Range("L99").Value = 2283.3333
$ws.Range("M99").Value = 2.333399999999983
$ws.Range("N99").Value = -5279.3333

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0

$ws.Range("H126").Value = 1692.5834
$ws.Range("I126").Value = 1495.6666
$ws.Range("J126").Value = 2283.3333
$ws.Range("K126").Value = 4486.9998
$ws.Range("L126").Value = 6849.999899999999
$ws.Range("M126").Value = -2016.9998
$ws.Range("N126").Value = -11789.9999

$ws.Range("H133").Value = 33333.168
$ws.Range("J133").Value = 30000
$ws.Range("L133").Value = 30000
$ws.Range("N133").Value = -35060

$ws.Range("H134").Value = 2383145.2
$ws.Range("I134").Value = 2436.4443
$ws.Range("J134").Value = 23809524
$ws.Range("K134").Value = 7309.3329
$ws.Range("L134").Value = 71428572
$ws.Range("M134").Value = -4774.3329
$ws.Range("N134").Value = -71433642

$ws.Range("H141").Value = 36160.25
$ws.Range("I141").Value = 15000
$ws.Range("J141").Value = 39183.145
$ws.Range("K141").Value = 15000
$ws.Range("L141").Value = 39183.145
$ws.Range("M141").Value = -9820
$ws.Range("N141").Value = -49543.145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1800
$ws.Range("J75").Value = 1800
$ws.Range("L75").Value = 5400
$ws.Range("N75").Value = -7396

$ws.Range("H78").Value = 1800
$ws.Range("J78").Value = 1800
$ws.Range("L78").Value = 16200
$ws.Range("N78").Value = -26184

$ws.Range("H107").Value = 406.79486
$ws.Range("I107").Value = 188.92857
$ws.Range("J107").Value = 528.8
$ws.Range("K107").Value = 566.78571
$ws.Range("L107").Value = 1586.4
$ws.Range("M107").Value = 1353.21429
$ws.Range("N107").Value = -5426.4

$ws.Range("H132").Value = 1508.88
$ws.Range("I132").Value = 706.3333
$ws.Range("J132").Value = 1762.3158
$ws.Range("K132").Value = 6356.9997
$ws.Range("L132").Value = 15860.8422
$ws.Range("M132").Value = -3826.9997
$ws.Range("N132").Value = -20920.8422

$ws.Range("H140").Value = 2856.05
$ws.Range("I140").Value = 1735.0834
$ws.Range("J140").Value = 4537.5
$ws.Range("K140").Value = 5205.2502
$ws.Range("L140").Value = 13612.5
$ws.Range("M140").Value = -25.2502000000004
$ws.Range("N140").Value = -23972.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2694.4736
$ws.Range("I80").Value = 2659.6
$ws.Range("J80").Value = 2706.9285
$ws.Range("K80").Value = 2659.6
$ws.Range("L80").Value = 2706.9285
$ws.Range("M80").Value = -1661.6
$ws.Range("N80").Value = -4702.9285

$ws.Range("H83").Value = 2694.4736
$ws.Range("I83").Value = 2659.6
$ws.Range("J83").Value = 2706.9285
$ws.Range("K83").Value = 13298
$ws.Range("L83").Value = 13534.6425
$ws.Range("M83").Value = -8306
$ws.Range("N83").Value = -23518.6425

$ws.Range("H122").Value = 2288.6191
$ws.Range("I122").Value = 2130.7334
$ws.Range("J122").Value = 2683.3333
$ws.Range("K122").Value = 6392.2002
$ws.Range("L122").Value = 8049.999899999999
$ws.Range("M122").Value = -3942.2002
$ws.Range("N122").Value = -12949.9999

$ws.Range("H126").Value = 4037.1875
$ws.Range("I126").Value = 2759.6
$ws.Range("J126").Value = 4617.909
$ws.Range("K126").Value = 8278.799999999999
$ws.Range("L126").Value = 13853.727
$ws.Range("M126").Value = -5808.799999999999
$ws.Range("N126").Value = -18793.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1139.862
$ws.Range("I22").Value = 923.3333
$ws.Range("J22").Value = 1164.8462
$ws.Range("K22").Value = 923.3333
$ws.Range("L22").Value = 1164.8462
$ws.Range("M22").Value = -628.3333
$ws.Range("N22").Value = -1754.8462

$ws.Range("H27").Value = 1139.862
$ws.Range("I27").Value = 923.3333
$ws.Range("J27").Value = 1164.8462
$ws.Range("K27").Value = 923.3333
$ws.Range("L27").Value = 1164.8462
$ws.Range("M27").Value = -816.3333
$ws.Range("N27").Value = -1378.8462

$ws.Range("H40").Value = 6183.3335
$ws.Range("I40").Value = 6680
$ws.Range("J40").Value = 3700
$ws.Range("K40").Value = 6680
$ws.Range("L40").Value = 3700
$ws.Range("M40").Value = -6544
$ws.Range("N40").Value = -3972

$ws.Range("H55").Value = 418.8421
$ws.Range("I55").Value = 261.1111
$ws.Range("J55").Value = 560.8
$ws.Range("K55").Value = 261.1111
$ws.Range("L55").Value = 560.8
$ws.Range("M55").Value = -88.11110000000002
$ws.Range("N55").Value = -906.8

$ws.Range("H122").Value = 6358.7407
$ws.Range("I122").Value = 5307.4287
$ws.Range("J122").Value = 7490.923
$ws.Range("K122").Value = 15922.2861
$ws.Range("L122").Value = 22472.769
$ws.Range("M122").Value = -13472.2861
$ws.Range("N122").Value = -27372.769

$ws.Range("H132").Value = 20008106
$ws.Range("I132").Value = 3973.4443
$ws.Range("J132").Value = 71447304
$ws.Range("K132").Value = 11920.3329
$ws.Range("L132").Value = 214341912
$ws.Range("M132").Value = -9390.332900000001
$ws.Range("N132").Value = -214346972

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6272.909
$ws.Range("I62").Value = 5233.6665
$ws.Range("J62").Value = 7520
$ws.Range("K62").Value = 5233.6665
$ws.Range("L62").Value = 7520
$ws.Range("M62").Value = -4609.6665
$ws.Range("N62").Value = -8768

$ws.Range("H65").Value = 6272.909
$ws.Range("I65").Value = 5233.6665
$ws.Range("J65").Value = 7520
$ws.Range("K65").Value = 26168.3325
$ws.Range("L65").Value = 37600
$ws.Range("M65").Value = -23048.3325
$ws.Range("N65").Value = -43840
